$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# Insert 5 new rows starting at row 95 (shifts old rows 95-110 down to 100-115)
$ws.Rows.Item(95).Resize(5).Insert()

# Fill the 5 newly inserted rows (95-99) with the new Artisan command documentation
$ws.Range("B95").Value = "keyboard(<bool>)"
$ws.Range("C95").Value = "enables/disables keyboard mode"

$ws.Range("B96").Value = "showCurve(<name>,<bool>)"
$ws.Range("C96").Value = "shows/hides the curve indicated by <name> which is one of { ET, BT, DeltaET, DeltaBT, BackgroundET, BackgroundBT}"

$ws.Range("B97").Value = "showExtraCurve(<extra_device>,<curve>,<bool>)"
$ws.Range("C97").Value = "shows/hides the <curve> (one of {T1,T2}) of the zero-based <extra_device> number"

$ws.Range("B98").Value = "showEvents(<event_type>, <bool>)"
$ws.Range("C98").Value = "shows/hides the events of <event_type> in [1,..,5]"

$ws.Range("B99").Value = "showBackgroundEvents(<bool>)"
$ws.Range("C99").Value = "shows/hides the events of the background profile"

# Set row heights to match the target layout (rows 100-111 already carry the
# default row height of 15 after the insert/shift, so no explicit change needed there)
$ws.Rows.Item(95).RowHeight = 13.8
$ws.Rows.Item(96).RowHeight = 13.8
$ws.Rows.Item(97).RowHeight = 13.8
$ws.Rows.Item(98).RowHeight = 13.8
$ws.Rows.Item(99).RowHeight = 13.8

# Update the view state: selection on C97, scrolled so row 85 is at the top
$ws.Activate()
$ws.Range("C97").Select()
$excel.ActiveWindow.ScrollRow = 85
$excel.ActiveWindow.ScrollColumn = 1
